$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 44.946245505921233
$ws.Range("C2").Value = -8.6796474851898218
$ws.Range("D2").Value = -2.6724587353141138
$ws.Range("E2").Value = 5.1667374473021255

$ws.Range("B3").Value = 42.549059593217265
$ws.Range("C3").Value = 7.1691463110793299
$ws.Range("D3").Value = -18.927945376511605
$ws.Range("E3").Value = 30.181396220663299

$ws.Range("B1:E3").Select()
